$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.07611394332627118,
    0.01491988877118644,
    0.03820511122881356,
    0.01813923463983051,
    0.02425185381355932,
    0.06498113082627119,
    0.008598715572033898,
    0.01576403601694915,
    0.05594213453389831,
    0.00998245497881356,
    0.02977356991525424,
    0.04193425582627119,
    0.02400357521186441,
    0.05972093485169491,
    0.09224874205508475,
    0.02147113347457627,
    0.03576370497881356,
    0.04141452595338983,
    0.04653237552966102,
    0.02200907044491525,
    0.001931607521186441,
    0.007435116525423729,
    0.003454382944915254,
    0.01107653601694915,
    0.007027939618644068,
    0.002219610699152542,
    0.0001870365466101695,
    0.01541313559322034,
    0.01811275158898305,
    0.002484441207627119,
    0.005877582097457627,
    0.02192300052966102,
    0.1603515625
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
